$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$rng = $ws.Range("A20:H20")
$rng.Interior.ThemeColor = 10
$rng.Interior.TintAndShade = 0
Write-Host "set ok"
